$wb = $excel.ActiveWorkbook
$sys = $wb.Worksheets.Item("#system")

$sys.Range("G68").Value2 = "saveModalDialogTextByLocator(var,locator)"
$sys.Range("G85").Value2 = "typeKeys(os,keystrokes)"
$sys.Range("G86").Value2 = "typeTextArea(name,text1,text2,text3,text4)"
$sys.Range("G87").Value2 = "typeTextBox(name,text1,text2,text3,text4)"
$sys.Range("G88").Value2 = "useApp(appId)"
$sys.Range("G89").Value2 = "useForm(formName)"
$sys.Range("G90").Value2 = "useHierTable(var,name)"
$sys.Range("G91").Value2 = "useList(var,name)"
$sys.Range("G92").Value2 = "useTable(var,name)"
$sys.Range("G93").Value2 = "useTableRow(var,row)"
$sys.Range("G94").Value2 = "waitFor(name,maxWaitMs)"
$sys.Range("G95").Value2 = "waitForLocator(locator,maxWaitMs)"
$sys.Range("J7").Value2 = "saveDiff(var,baseline,actual)"
$sys.Range("M12").Value2 = "compact(var,json,removeEmpty)"
$sys.Range("M13").Value2 = "fromCsv(csv,header,jsonFile)"
$sys.Range("M14").Value2 = "minify(json,var)"
$sys.Range("M15").Value2 = "storeCount(json,jsonpath,var)"
$sys.Range("M16").Value2 = "storeValue(json,jsonpath,var)"
$sys.Range("M17").Value2 = "storeValues(json,jsonpath,var)"
$sys.Range("Y50").Value2 = "clickAll(locator)"
$sys.Range("Y51").Value2 = "clickAndWait(locator,waitMs)"
$sys.Range("Y52").Value2 = "clickByLabel(label)"
$sys.Range("Y53").Value2 = "clickByLabelAndWait(label,waitMs)"
$sys.Range("Y54").Value2 = "clickOffset(locator,x,y)"
$sys.Range("Y55").Value2 = "clickWithKeys(locator,keys)"
$sys.Range("Y56").Value2 = "close()"
$sys.Range("Y57").Value2 = "closeAll()"
$sys.Range("Y58").Value2 = "deselect(locator,text)"
$sys.Range("Y59").Value2 = "deselectMulti(locator,array)"
$sys.Range("Y60").Value2 = "dismissInvalidCert()"
$sys.Range("Y61").Value2 = "dismissInvalidCertPopup()"
$sys.Range("Y62").Value2 = "doubleClick(locator)"
$sys.Range("Y63").Value2 = "doubleClickAndWait(locator,waitMs)"
$sys.Range("Y64").Value2 = "doubleClickByLabel(label)"
$sys.Range("Y65").Value2 = "doubleClickByLabelAndWait(label,waitMs)"
$sys.Range("Y66").Value2 = "dragAndDrop(fromLocator,toLocator)"
$sys.Range("Y67").Value2 = "dragTo(fromLocator,xOffset,yOffset)"
$sys.Range("Y68").Value2 = "editLocalStorage(key,value)"
$sys.Range("Y69").Value2 = "executeScript(var,script)"
$sys.Range("Y70").Value2 = "focus(locator)"
$sys.Range("Y71").Value2 = "goBack()"
$sys.Range("Y72").Value2 = "goBackAndWait()"
$sys.Range("Y73").Value2 = "maximizeWindow()"
$sys.Range("Y74").Value2 = "mouseOver(locator)"
$sys.Range("Y75").Value2 = "open(url)"
$sys.Range("Y76").Value2 = "openAndWait(url,waitMs)"
$sys.Range("Y77").Value2 = "openHttpBasic(url,username,password)"
$sys.Range("Y78").Value2 = "openIgnoreTimeout(url)"
$sys.Range("Y79").Value2 = "refresh()"
$sys.Range("Y80").Value2 = "refreshAndWait()"
$sys.Range("Y81").Value2 = "resizeWindow(width,height)"
$sys.Range("Y82").Value2 = "rightClick(locator)"
$sys.Range("Y83").Value2 = "saveAllWindowIds(var)"
$sys.Range("Y84").Value2 = "saveAllWindowNames(var)"
$sys.Range("Y85").Value2 = "saveAttribute(var,locator,attrName)"
$sys.Range("Y86").Value2 = "saveAttributeList(var,locator,attrName)"
$sys.Range("Y87").Value2 = "saveCount(var,locator)"
$sys.Range("Y88").Value2 = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$sys.Range("Y89").Value2 = "saveElement(var,locator)"
$sys.Range("Y90").Value2 = "saveElements(var,locator)"
$sys.Range("Y91").Value2 = "saveLocalStorage(var,key)"
$sys.Range("Y92").Value2 = "saveLocation(var)"
$sys.Range("Y93").Value2 = "savePageAs(var,sessionIdName,url)"
$sys.Range("Y94").Value2 = "savePageAsFile(sessionIdName,url,file)"
$sys.Range("Y95").Value2 = "saveTableAsCsv(locator,nextPageLocator,file)"
$sys.Range("Y96").Value2 = "saveText(var,locator)"
$sys.Range("Y97").Value2 = "saveTextArray(var,locator)"
$sys.Range("Y98").Value2 = "saveTextSubstringAfter(var,locator,delim)"
$sys.Range("Y99").Value2 = "saveTextSubstringBefore(var,locator,delim)"
$sys.Range("Y100").Value2 = "saveTextSubstringBetween(var,locator,start,end)"
$sys.Range("Y101").Value2 = "saveValue(var,locator)"
$sys.Range("Y102").Value2 = "saveValues(var,locator)"
$sys.Range("Y103").Value2 = "scrollElement(locator,xOffset,yOffset)"
$sys.Range("Y104").Value2 = "scrollLeft(locator,pixel)"
$sys.Range("Y105").Value2 = "scrollPage(xOffset,yOffset)"
$sys.Range("Y106").Value2 = "scrollRight(locator,pixel)"
$sys.Range("Y107").Value2 = "scrollTo(locator)"
$sys.Range("Y108").Value2 = "select(locator,text)"
$sys.Range("Y109").Value2 = "selectFrame(locator)"
$sys.Range("Y110").Value2 = "selectMulti(locator,array)"
$sys.Range("Y111").Value2 = "selectMultiOptions(locator)"
$sys.Range("Y112").Value2 = "selectText(locator)"
$sys.Range("Y113").Value2 = "selectWindow(winId)"
$sys.Range("Y114").Value2 = "selectWindowAndWait(winId,waitMs)"
$sys.Range("Y115").Value2 = "selectWindowByIndex(index)"
$sys.Range("Y116").Value2 = "selectWindowByIndexAndWait(index,waitMs)"
$sys.Range("Y117").Value2 = "toggleSelections(locator)"
$sys.Range("Y118").Value2 = "type(locator,value)"
$sys.Range("Y119").Value2 = "typeKeys(locator,value)"
$sys.Range("Y120").Value2 = "uncheckAll(locator)"
$sys.Range("Y121").Value2 = "unselectAllText()"
$sys.Range("Y122").Value2 = "upload(fieldLocator,file)"
$sys.Range("Y123").Value2 = "verifyContainText(locator,text)"
$sys.Range("Y124").Value2 = "verifyText(locator,text)"
$sys.Range("Y125").Value2 = "wait(waitMs)"
$sys.Range("Y126").Value2 = "waitForElementPresent(locator)"
$sys.Range("Y127").Value2 = "waitForPopUp(winId,waitMs)"
$sys.Range("Y128").Value2 = "waitForTextPresent(text)"
$sys.Range("Y129").Value2 = "waitForTitle(text)"
$sys.Range("AD7").Value2 = "assertSoap(wsdl,xml)"
$sys.Range("AD8").Value2 = "assertSoapFaultCode(expected,xml)"
$sys.Range("AD9").Value2 = "assertSoapFaultString(expected,xml)"
$sys.Range("AD10").Value2 = "assertValue(xml,xpath,expected)"
$sys.Range("AD11").Value2 = "assertValues(xml,xpath,array,exactOrder)"
$sys.Range("AD12").Value2 = "assertWellformed(xml)"
$sys.Range("AD13").Value2 = "beautify(xml,var)"
$sys.Range("AD14").Value2 = "clear(xml,xpath,var)"
$sys.Range("AD15").Value2 = "delete(xml,xpath,var)"
$sys.Range("AD16").Value2 = "insertAfter(xml,xpath,content,var)"
$sys.Range("AD17").Value2 = "insertBefore(xml,xpath,content,var)"
$sys.Range("AD18").Value2 = "minify(xml,var)"
$sys.Range("AD19").Value2 = "prepend(xml,xpath,content,var)"
$sys.Range("AD20").Value2 = "replace(xml,xpath,content,var)"
$sys.Range("AD21").Value2 = "replaceIn(xml,xpath,content,var)"
$sys.Range("AD22").Value2 = "storeCount(xml,xpath,var)"
$sys.Range("AD23").Value2 = "storeSoapFaultCode(var,xml)"
$sys.Range("AD24").Value2 = "storeSoapFaultDetail(var,xml)"
$sys.Range("AD25").Value2 = "storeSoapFaultString(var,xml)"
$sys.Range("AD26").Value2 = "storeValue(xml,xpath,var)"
$sys.Range("AD27").Value2 = "storeValues(xml,xpath,var)"

$names = $wb.Names
foreach ($n in $names) {
    if ($n.Name -eq "desktop") { $n.RefersTo = "='#system'!`$G`$2:`$G`$95" }
    if ($n.Name -eq "image") { $n.RefersTo = "='#system'!`$J`$2:`$J`$7" }
    if ($n.Name -eq "json") { $n.RefersTo = "='#system'!`$M`$2:`$M`$17" }
    if ($n.Name -eq "web") { $n.RefersTo = "='#system'!`$Y`$2:`$Y`$129" }
    if ($n.Name -eq "xml") { $n.RefersTo = "='#system'!`$AD`$2:`$AD`$27" }
}
